$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.83935558795929
$ws.Range("B1").Value = 4.856363296508789
$ws.Range("C1").Value = 3.686752319335938
$ws.Range("D1").Value = 1.251519918441772
$ws.Range("E1").Value = 0.8272326588630676
